{"js": "// Product line-item row and its \"discounted price\" sibling row (the 3rd and\n// 4th rows, 0-based index 2 and 3, of the 4th table in the document body)\n// get a shorter row height and several cells re-centered.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[3];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst targetRowIndexes = [2, 3];\nconst rows = targetRowIndexes.map((i) => table.rows.items[i]);\n\n// Shrink the row height: 432 twips (21.6pt) -> 288 twips (14.4pt).\nfor (const row of rows) {\n  row.preferredHeight = 14.4;\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows) {\n  for (let col = 0; col < 5; col++) {\n    row.cells.items[col].body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\n// Columns 0, 1, 3, 4 become centered; column 2 drops its explicit\n// \"center\" alignment back to the implicit default (left).\nfor (const row of rows) {\n  for (let col = 0; col < 5; col++) {\n    const paragraph = row.cells.items[col].body.paragraphs.items[0];\n    paragraph.alignment = col === 2 ? \"Left\" : \"Centered\";\n  }\n}\nawait context.sync();\n", "ps1": "# Product line-item row and its \"discounted price\" sibling row (rows 3 and 4,\n# 1-based, of the 4th table in the document) get a shorter row height and\n# several cells re-centered.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(4)\n\nforeach ($r in 3, 4) {\n    $row = $tbl.Rows.Item($r)\n\n    # Shrink the row height: 432 twips (21.6pt) -> 288 twips (14.4pt).\n    $row.Height = 14.4\n\n    # Columns 1, 2, 4, 5 become centered; column 3 drops its explicit\n    # \"center\" alignment back to the implicit default (left).\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $paragraph = $cell.Range.Paragraphs.Item(1)\n        if ($c -eq 3) {\n            $paragraph.Alignment = 0\n        } else {\n            $paragraph.Alignment = 1\n        }\n    }\n}\n"}
